{"js": "// Load all body paragraphs so we can locate the title heading and the\n// two paragraphs at the very end of the document.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- 1. Insert a new \"Meta description\" paragraph right after the H1 title ---\nconst titlePara = items[0]; // \"Play Chilli Master free: Game Review & Features\" (Heading1)\n\nconst metaPara = titlePara.insertParagraph(\"\", \"After\");\nmetaPara.style = \"Normal\"; // plain body paragraph (no heading style)\n\n// Insert the non-bold remainder first ...\nmetaPara.insertText(\n  \": Read our review of Chilli Master, a slot game with colorful graphics, cuisine-inspired symbols, and high volatility. Play for free today!\",\n  \"End\"\n);\nawait context.sync();\n\n// ... then insert the bold \"Meta description\" label at the start so its\n// own formatting doesn't bleed into the text inserted afterwards.\nconst metaLabelRange = metaPara.insertText(\"Meta description\", \"Start\");\nmetaLabelRange.font.bold = true;\nawait context.sync();\n\n// --- 2. Remove the duplicated bold title paragraph & update the italic one ---\nconst n = items.length;\nconst italicPara = items[n - 1]; // italic \"Read our review of Chilli Master...\" paragraph\nconst boldDuplicatePara = items[n - 2]; // bold \"Play Chilli Master free...\" paragraph\n\nboldDuplicatePara.delete();\nawait context.sync();\n\nitalicPara.insertText(\n  \"Create a cartoon-style feature image for \\\"Chilli Master\\\" that showcases a happy Maya warrior with glasses. The image should have a colorful and festive feel, with elements of Mexican cuisine and hot peppers incorporated into the design. Make sure the image captures the excitement and fun of the game, and reflects the game's theme of Mexican culture and cuisine.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Insert a new \"Meta description\" paragraph right after the H1 title ---\n$titlePara = $d.Paragraphs.Item(1)   # \"Play Chilli Master free: Game Review & Features\" (Heading1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs.Item(2)\n$metaPara.Style = $d.Styles.Item(\"Normal\")   # plain body paragraph (no heading style)\n\n$metaRange = $metaPara.Range\n$metaRange.Text = \": Read our review of Chilli Master, a slot game with colorful graphics, cuisine-inspired symbols, and high volatility. Play for free today!\"\n\n# Add the bold \"Meta description\" label in front of the text above.\n$labelInsertPoint = $d.Range($metaRange.Start, $metaRange.Start)\n$labelInsertPoint.InsertBefore(\"Meta description\")\n$labelRange = $d.Range($metaRange.Start, $metaRange.Start + (\"Meta description\").Length)\n$labelRange.Bold = 1\n\n# --- 2. Remove the duplicated bold title paragraph & update the italic one ---\n$count = $d.Paragraphs.Count\n$boldDuplicatePara = $d.Paragraphs.Item($count - 1)   # bold \"Play Chilli Master free...\" paragraph\n$italicPara = $d.Paragraphs.Item($count)              # italic \"Read our review of Chilli Master...\" paragraph\n\n$italicRange = $italicPara.Range\n$oldLen = $italicRange.Text.Length\n$oldStart = $italicRange.Start\n$newText = \"Create a cartoon-style feature image for `\"Chilli Master`\" that showcases a happy Maya warrior with glasses. The image should have a colorful and festive feel, with elements of Mexican cuisine and hot peppers incorporated into the design. Make sure the image captures the excitement and fun of the game, and reflects the game's theme of Mexican culture and cuisine.\"\n\n# Insert the replacement text before the old text (preserves exact characters,\n# no smart-quote autocorrection), then explicitly re-apply italic formatting,\n# then delete the now-stale old text that follows.\n$newTextInsertPoint = $d.Range($oldStart, $oldStart)\n$newTextInsertPoint.InsertBefore($newText)\n$newTextRange = $d.Range($oldStart, $oldStart + $newText.Length)\n$newTextRange.Italic = 1\n\n$oldTextRange = $d.Range($oldStart + $newText.Length, $oldStart + $newText.Length + $oldLen)\n$oldTextRange.Delete()\n\n# Finally, remove the whole duplicated bold title paragraph.\n$boldDuplicatePara.Range.Delete()\n"}
